$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New timestamp that applies to every data row after this run.
$newTimestamp = "2025-12-26 18:27:04"

# --- Insert two fresh rows above the existing data (rows 2 & 3 push down to 4 & 5) ---
$ws.Rows("2:3").Insert()

# --- Column width tweaks (B: 30 -> 41, H: 12 -> 13) ---
# ColumnWidth is in "character" units and Excel stores the serialized width
# with a fixed +0.8333... padding, so compensate to land exactly on 41 / 13.
$ws.Columns("B").ColumnWidth = 40.166666666666664
$ws.Columns("H").ColumnWidth = 12.166666666666666

# --- Row 2 (new): AI tech-lead job post ---
$ws.Range("A2").Value = $newTimestamp
$ws.Range("B2").Value = "大企業の業務効率化AIプロジェクトの技術方針策定を支援するAIテックリード募集"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5423720"
$ws.Range("G2").Value = 385
$ws.Range("H2").Value = "🔥AI,Ai ◆効率化"

# --- Row 3 (new): AI 1-hour online teaching job post ---
$ws.Range("A3").Value = $newTimestamp
$ws.Range("B3").Value = "AIについて1時間オンラインで教えて欲しい"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5461891"
$ws.Range("G3").Value = 295
$ws.Range("H3").Value = "🔥AI,Ai"

# --- Rows 4 & 5 (previously 2 & 3): only the captured timestamp changes ---
$ws.Range("A4").Value = $newTimestamp
$ws.Range("A5").Value = $newTimestamp

# --- Rebuild hyperlinks for F2:F5 (stale refs from the row insert are cleared first) ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5423720")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5461891")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5461280")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5461481")

# Hyperlinks.Add leaves the cell's own style untouched on this engine for F4/F5
# (already "Hyperlink" from the row shift) but F2/F3 are brand-new cells, so
# stamp them with the same visual style used by the other link cells.
$ws.Range("F2").Style = "Hyperlink"
$ws.Range("F3").Style = "Hyperlink"
$ws.Range("F4").Style = "Hyperlink"
$ws.Range("F5").Style = "Hyperlink"
